$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.817.09"
$ws.Range("E2").Value = "  +2.58%  "
$ws.Range("D3").Value = "1.767.34"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("E4").Value = "  -0.51%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.99"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9963"
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3818"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3424"
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.02"
$ws.Range("E9").Value = "  -2.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.138"
$ws.Range("E10").Value = "  -4.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07403"
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9975"
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.37"
$ws.Range("E13").Value = "  +3.54%  "
$ws.Range("E14").Value = "  -1.07%  "
$ws.Range("D15").Value = "1.768.19"
$ws.Range("E15").Value = "  -0.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.089"
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001074"
$ws.Range("E17").Value = "  -0.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06662"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "82.21"
$ws.Range("E19").Value = "  -1.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9976"
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.35"
$ws.Range("E21").Value = "  +0.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.416"
$ws.Range("E22").Value = "  -1.97%  "
$ws.Range("D23").Value = "27.842.14"
$ws.Range("E23").Value = "  +2.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.04"
$ws.Range("E24").Value = "  -1.30%  "
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.440"
$ws.Range("E26").Value = "  -1.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.71"
$ws.Range("E27").Value = "  -1.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.415"
$ws.Range("E28").Value = "  -3.59%  "
$ws.Range("E29").Value = "  -0.95%  "
$ws.Range("B30").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C30").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D30").Value = "1.969.84"
$ws.Range("E30").Value = "  -0.47%  "
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "134.22"
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.140"
$ws.Range("E32").Value = "  +2.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.951"
$ws.Range("E33").Value = "  -1.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08782"
$ws.Range("E34").Value = "  +1.53%  "
$ws.Range("E35").Value = "  -2.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02435"
$ws.Range("E36").Value = "  +5.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6837"
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.312"
$ws.Range("E38").Value = "  -1.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06327"
$ws.Range("E39").Value = "  +0.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2185"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.515"
$ws.Range("E41").Value = "  -6.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.233"
$ws.Range("E42").Value = "  +0.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.242"
$ws.Range("E43").Value = "  -3.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.19"
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9962"
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6268"
$ws.Range("E46").Value = "  -2.04%  "
$ws.Range("E47").Value = "  -0.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.49"
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.079"
$ws.Range("E49").Value = "  -2.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07381"
$ws.Range("E50").Value = "  +4.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.25"
$ws.Range("E51").Value = "  -0.39%  "
